$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the "id" column (A). Everything shifts left by one. ---
$ws.Columns.Item(1).Delete()

# --- 2. Fix up the tail of the header row. ---
# After the shift, AF1/AG1 currently read "employment_status"/"age_group";
# the edited template swaps them (age_group first, with a trailing space) and
# un-bolds them, and the four trailing header cells (old interview_date_time,
# interview_status, last_downloaded_date + the now-empty last slot) are wiped out.
$ws.Range("AF1").Value = "age_group "
$ws.Range("AG1").Value = "employment_status"
$ws.Range("AF1:AG1").Font.Bold = $false
$ws.Range("AH1:AJ1").Clear()

# --- 3. Row 1 formatting: slightly taller, and A1 now carries the bold header style too. ---
$ws.Rows.Item(1).RowHeight = 14.9
$ws.Range("A1").Font.Bold = $true

# --- 4. Respondent / data fixes (columns already shifted left by one). ---
$ws.Range("D2").Value = "Respondent Seven"
$ws.Range("E2").Value = 73234575

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Respondent Eight"
$ws.Range("E3").Value = 73345665

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "Respondent Nine"
$ws.Range("E4").Value = 73356755

# D3/D4 pick up their own (non-bold) style, distinct from the plain default cells.
$ws.Range("D3:D4").Font.Bold = $false

# --- 5. Column width touch-ups to match the edited template. ---
$ws.Columns.Item(2).ColumnWidth = 10.47
$ws.Columns.Item(22).ColumnWidth = 11.24
$ws.Columns.Item(23).ColumnWidth = 10.14
$ws.Columns.Item(26).ColumnWidth = 21.39
$ws.Columns.Item(28).ColumnWidth = 16.1
$ws.Columns.Item(30).ColumnWidth = 8.05
$ws.Columns.Item(32).ColumnWidth = 11.68
$ws.Columns.Item(33).ColumnWidth = 21.28

# --- 6. Stray formatting far out on row 1 (mirrors the legacy AMG:AMJ artifact). ---
$ws.Range("AMG1:AMJ1").NumberFormat = "General"

# --- 7. Sheet view tweaks. ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E4").Select()
